# Applies the "invalidFor.docx" fixture edit:
#   - keeps the existing "+++FOR company IN companies+++" / "+++INS $company.name+++"
#     paragraphs untouched
#   - turns the old (invalid) "+++END-FOR invalidKey+++" paragraph into a second,
#     nested FOR loop over `persons` (with its own INS) followed by the matching
#     END-FOR for `company` (re-using the existing `_GoBack` bookmark), with blank
#     separator paragraphs around the loops, matching the target OOXML.

$d = $word.ActiveDocument
$wns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

function New-Run($text, [switch]$preserve) {
    if ($preserve) {
        return "<w:r><w:t xml:space='preserve'>" + $text + "</w:t></w:r>"
    }
    return "<w:r><w:t>" + $text + "</w:t></w:r>"
}

# ---------------------------------------------------------------------------
# Step 1: insert a brand-new empty paragraph right before the (old) third
# paragraph ("+++END-FOR invalidKey+++"). This becomes the blank separator
# paragraph between the `company` INS and the new `person` FOR.
# ---------------------------------------------------------------------------
$p3 = $d.Paragraphs.Item(3)
$ip = $d.Range($p3.Range.Start, $p3.Range.Start)
$null = $ip.InsertXML("<w:p $wns></w:p>")

# ---------------------------------------------------------------------------
# Step 2: the old "+++END-FOR invalidKey+++" paragraph (with its proofErr-
# wrapped run) is now paragraph 4. Replace its whole contents with the new
# "+++FOR person IN persons+++" run sequence.
# ---------------------------------------------------------------------------
$p4 = $d.Paragraphs.Item(4)
$target = $d.Range($p4.Range.Start, $p4.Range.End)
$xml = "<w:p $wns>" `
    + (New-Run "+++FOR " -preserve) `
    + (New-Run "person") `
    + (New-Run " IN " -preserve) `
    + (New-Run "persons") `
    + (New-Run "+++") `
    + "</w:p>"
$null = $target.InsertXML($xml)

# ---------------------------------------------------------------------------
# Step 3: insert the "+++INS $person.name+++" paragraph right after it
# (before the bookmark paragraph, which is still paragraph 5 at this point).
# ---------------------------------------------------------------------------
$p5 = $d.Paragraphs.Item(5)
$ip2 = $d.Range($p5.Range.Start, $p5.Range.Start)
$xml = "<w:p $wns>" `
    + (New-Run "+++INS `$") `
    + (New-Run "person") `
    + (New-Run ".name+++") `
    + "</w:p><w:p $wns></w:p>"
$null = $ip2.InsertXML($xml)

# ---------------------------------------------------------------------------
# Step 4: the bookmark-only paragraph (formerly paragraph 4) is now
# paragraph 6. Replace its contents, keeping the existing `_GoBack`
# bookmark in place but adding "+++END-FOR " before it, "company" between
# bookmarkStart/bookmarkEnd, and "+++" after it.
# ---------------------------------------------------------------------------
$p6 = $d.Paragraphs.Item(6)
$target2 = $d.Range($p6.Range.Start, $p6.Range.End)
$xml = "<w:p $wns>" `
    + (New-Run "+++END-FOR " -preserve) `
    + "<w:bookmarkStart w:id='0' w:name='_GoBack'/>" `
    + (New-Run "company") `
    + "<w:bookmarkEnd w:id='0'/>" `
    + (New-Run "+++") `
    + "</w:p>"
$null = $target2.InsertXML($xml)

# ---------------------------------------------------------------------------
# Step 5: after that paragraph, add a blank separator paragraph followed by
# the final "+++END-FOR company+++" paragraph, and one last trailing blank
# paragraph.
# ---------------------------------------------------------------------------
$p6b = $d.Paragraphs.Item(6)
$ip3 = $d.Range($p6b.Range.End, $p6b.Range.End)
$xml = "<w:p $wns></w:p><w:p $wns>" `
    + (New-Run "+++END-FOR " -preserve) `
    + (New-Run "company") `
    + (New-Run "+++") `
    + "</w:p><w:p $wns></w:p>"
$null = $ip3.InsertXML($xml)

Write-Host "Final paragraph count:" $d.Paragraphs.Count
foreach ($p in $d.Paragraphs) {
    Write-Host ("  > '" + $p.Range.Text + "'")
}
